$d = $word.ActiveDocument

# --- 1. Close out the sentence at the end of the "VR INTEGRATION" section
#        with a separate "." run (keeps the same en-GB language formatting). ---
$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastParaIndex)
$lastRange = $lastPara.Range
$lastRange.InsertAfter(".")
$lastRange.LanguageID = "en-GB"

# --- 2. New paragraph: progress note about adding stereo shaders for both
#        video and image. ---
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs($d.Paragraphs.Count)
$newRange = $newPara.Range
$newRange.InsertAfter("OK added both video and image shader and screen for both, very spaghetti and hackish way but works for now. Also wth recheck your shaders assignment etc properly, this took me like hours to debug bruh.")

# --- 3. New paragraph: testing with VR, the editor error, and the
#        exported-build fixes / remaining issues. ---
$newPara = $d.Paragraphs($d.Paragraphs.Count)
$newPara.Range.InsertParagraphAfter()

$p3 = $d.Paragraphs($d.Paragraphs.Count)
$r3 = $p3.Range
$r3.InsertAfter("Now test with VR.")

$p3 = $d.Paragraphs($d.Paragraphs.Count)
$r3 = $p3.Range
$r3.InsertAfter(" Ok for some reason my godot run in editor don’t work for VR, getting ")

$p3 = $d.Paragraphs($d.Paragraphs.Count)
$r3 = $p3.Range
$r3.InsertAfter("OpenXR: No viewport was marked with use_xr, there is no rendered output!")

$p3 = $d.Paragraphs($d.Paragraphs.Count)
$r3 = $p3.Range
$r3.InsertAfter(" Error. But after exporting as .exe and running it works as expected, there was some bugs that is squashed, kinda weird how it suddenly appear and didn’t realised it before, ")

$p3 = $d.Paragraphs($d.Paragraphs.Count)
$r3 = $p3.Range
$r3.InsertAfter("mainly the file=false Boolean change when pressing file button multiple times, and still unsure bug about SBS_Screen duplicates showing due to autoplay emit, need further research, other than that, I need to change the thumbnail gen to be better cus rn its so bad.")

$p3 = $d.Paragraphs($d.Paragraphs.Count)
$r3 = $p3.Range
$r3.InsertAfter(" (Add more delay)")
